$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), copying the existing header
# format (bold font, border, centered alignment) from H1 so the new
# headers match the style of the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "IF"

# Fill in the data values for columns I and J, rows 2-15.
$iValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 4
    14 = 1
    15 = 3
}

$jValues = @{
    2  = 5
    3  = 5
    4  = 5
    5  = 5
    6  = 5
    7  = 4
    8  = 5
    9  = 5
    10 = 6
    11 = 3
    12 = 4
    13 = 6
    14 = 3
    15 = 4
}

foreach ($row in 2..15) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
